# Set line spacing to single (1.0) for every paragraph in the document,
# i.e. <w:spacing w:line="240" w:lineRule="auto"/> in each <w:pPr>.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 0   # wdLineSpaceSingle
}
